# The report's closing line "Москва, 2024" needs its year bumped to
# "Москва, 2025" (report date fix, see commit "fixed report in lr4 && lr5-7").
$d = $word.ActiveDocument

$find = $d.Content.Find
$found = $find.Execute(
    "Москва, 2024",  # FindText
    $true,           # MatchCase
    $false,          # MatchWholeWord
    $false,          # MatchWildcards
    $false,          # MatchSoundsLike
    $false,          # MatchAllWordForms
    $true,           # Forward
    1,               # Wrap (wdFindContinue)
    $false,          # Format
    "Москва, 2025",  # ReplaceWith
    2                # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find 'Москва, 2024' to update to 'Москва, 2025'"
}
